$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.830990791320801
$ws.Range("B1").Value = 3.577049732208252
$ws.Range("C1").Value = 1.932513952255249
$ws.Range("D1").Value = 1.539733648300171
$ws.Range("E1").Value = 1.416112661361694
